$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting existing rows 34-135 down to 35-136.
$ws.Rows("34:34").Insert()

# Populate the newly inserted row 34 with the new data point.
$ws.Range("A34").Value = 3
$ws.Range("B34").Value = "Femacal de La Calera"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = 44497
$ws.Range("E34").Value = 5
$ws.Range("F34").Value = 100112010
$ws.Range("G34").Value = "Achicoria"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 95
$ws.Range("K34").Value = 6000
$ws.Range("L34").Value = 6500
$ws.Range("M34").Value = 6263
$ws.Range("N34").Value = "$/caja 16 unidades"
$ws.Range("O34").Value = "Provincia de Quillota"
$ws.Range("P34").Value = 391
$ws.Range("Q34").Value = 16
$ws.Range("R34").Value = "Hortaliza"
